$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price rows need to be inserted at the top of the
# "Chirimoya" history (rows 152-153), pushing the existing rows
# (old 152..243) down by two (new 154..245).
$ws.Rows("152:153").Insert()

# New row 152: Especial quality, 80 units, 2800 $/kilo
$ws.Cells.Item(152, 1).Value  = 10
$ws.Cells.Item(152, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(152, 3).Value  = "La Araucanía"
$ws.Cells.Item(152, 4).Value  = 45233
$ws.Cells.Item(152, 5).Value  = 9
$ws.Cells.Item(152, 6).Value  = "Fruta"
$ws.Cells.Item(152, 7).Value  = 100107
$ws.Cells.Item(152, 8).Value  = "Otros"
$ws.Cells.Item(152, 9).Value  = 100107002
$ws.Cells.Item(152, 10).Value = "Chirimoya"
$ws.Cells.Item(152, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(152, 12).Value = "Especial"
$ws.Cells.Item(152, 13).Value = 80
$ws.Cells.Item(152, 14).Value = 2800
$ws.Cells.Item(152, 15).Value = 2800
$ws.Cells.Item(152, 16).Value = 2800
$ws.Cells.Item(152, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(152, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(152, 19).Value = 2800
$ws.Cells.Item(152, 20).Value = 1

# New row 153: Primera quality, 140 units, 2500 $/kilo
$ws.Cells.Item(153, 1).Value  = 10
$ws.Cells.Item(153, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(153, 3).Value  = "La Araucanía"
$ws.Cells.Item(153, 4).Value  = 45233
$ws.Cells.Item(153, 5).Value  = 9
$ws.Cells.Item(153, 6).Value  = "Fruta"
$ws.Cells.Item(153, 7).Value  = 100107
$ws.Cells.Item(153, 8).Value  = "Otros"
$ws.Cells.Item(153, 9).Value  = 100107002
$ws.Cells.Item(153, 10).Value = "Chirimoya"
$ws.Cells.Item(153, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(153, 12).Value = "Primera"
$ws.Cells.Item(153, 13).Value = 140
$ws.Cells.Item(153, 14).Value = 2500
$ws.Cells.Item(153, 15).Value = 2500
$ws.Cells.Item(153, 16).Value = 2500
$ws.Cells.Item(153, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(153, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(153, 19).Value = 2500
$ws.Cells.Item(153, 20).Value = 1
